$d = $word.ActiveDocument

# Paragraph: "{m:('dh' + i + '.gif').asImage().setWidth(100)}"
# Split run "{m" into "{" and "m" (split after the opening brace).
$r = $d.Range(73, 74)
$r.Delete()
$ins = $d.Range(73, 73)
$ins.InsertBefore("{")

# Split run ".setWidth(100)}" into ".setWidth(100)" and "}" (split before the closing brace).
$r = $d.Range(119, 120)
$r.Delete()
$left = $d.Range(105, 119)
$left.InsertAfter("}")

# Paragraph: "{m:endfor}"
# Split run "{m:" into "{" and "m:" (split after the opening brace).
$r = $d.Range(121, 122)
$r.Delete()
$ins = $d.Range(121, 121)
$ins.InsertBefore("{")
